$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '26.662.05'
$ws.Range('E2').Value = '  -1.72%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.591.44'
$ws.Range('E3').Value = '  -2.31%  '

$ws.Range('E4').Value = '  +0.05%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '211.40'
$ws.Range('E5').Value = '  -2.19%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.506'
$ws.Range('E6').Value = '  -1.68%  '

$ws.Range('E7').Value = '  +0.04%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.249'
$ws.Range('E8').Value = '  -1.50%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.0619'
$ws.Range('E9').Value = '  -0.64%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '19.63'
$ws.Range('E10').Value = '  -2.20%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0835'
$ws.Range('E11').Value = '  -1.88%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.814.58'
$ws.Range('E12').Value = '  -2.29%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.615.69'
$ws.Range('E13').Value = '  -1.09%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '4.05'
$ws.Range('E14').Value = '  -1.31%  '

$ws.Range('E15').Value = '  -1.71%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '63.86'
$ws.Range('E16').Value = '  -1.67%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '26.702.68'
$ws.Range('E17').Value = '  -1.43%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.0₃0731'
$ws.Range('E18').Value = '  -0.01%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '209.80'
$ws.Range('E19').Value = '  -1.82%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.70'
$ws.Range('E21').Value = '  -1.82%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.27'

$ws.Range('E23').Value = '  -4.13%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '8.89'
$ws.Range('E24').Value = '  -1.82%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '146.67'
$ws.Range('E25').Value = '  -0.53%  '

$ws.Range('E26').Value = '  +2.17%  '

$ws.Range('E27').Value = '  +0.10%  '

$ws.Range('E28').Value = '  -3.93%  '

$ws.Range('E29').Value = '  -1.02%  '

$ws.Range('E30').Value = '  -0.56%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.15'
$ws.Range('E31').Value = '  -1.87%  '

$ws.Range('E32').Value = '  -3.25%  '

$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.659'
$ws.Range('E33').Value = '  +22.96%  '

$ws.Range('B34').Value = 'InternetComputer(DFINITY)'
$ws.Range('C34').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '2.94'
$ws.Range('E34').Value = '  -1.89%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.309.71'
$ws.Range('E35').Value = '  -0.01%  '

$ws.Range('E36').Value = '  -3.17%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.42'
$ws.Range('E37').Value = '  -0.88%  '

$ws.Range('E38').Value = '  -0.74%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.821'
$ws.Range('E39').Value = '  -2.48%  '

$ws.Range('E40').Value = '  -0.01%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.787'
$ws.Range('E41').Value = '  -1.90%  '

$ws.Range('E42').Value = '  -4.38%  '

$ws.Range('E43').Value = '  +0.96%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '63.14'
$ws.Range('E44').Value = '  +0.92%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.728.66'
$ws.Range('E45').Value = '  -2.05%  '

$ws.Range('B46').Value = 'Quant'
$ws.Range('C46').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '89.01'
$ws.Range('E46').Value = '  -1.78%  '

$ws.Range('B47').Value = 'RenderToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.63'
$ws.Range('E47').Value = '  +1.87%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.824'
$ws.Range('E48').Value = '  +3.32%  '

$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0508'
$ws.Range('E49').Value = '  -0.96%  '

$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0981'
$ws.Range('E50').Value = '  +3.80%  '

$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '7.47'
$ws.Range('E51').Value = '  -0.79%  '
